# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# to reflect refreshed cryptos list data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.096.96"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "2.618.43"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'521.15"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").Value = "'148.61"
$ws.Range("E6").Value = "  -3.28%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'0.570"
$ws.Range("E8").Value = "  -4.74%  "
$ws.Range("D9").Value = "2.622.44"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("D10").Value = "'6.30"
$ws.Range("E10").Value = "  -5.43%  "
$ws.Range("D11").Value = "'0.105"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "'0.341"
$ws.Range("E12").Value = "  -1.71%  "
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").Value = "3.081.20"
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").Value = "60.101.99"
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").Value = "'21.16"
$ws.Range("E16").Value = "  -2.57%  "
$ws.Range("E17").Value = "  -1.77%  "
$ws.Range("D18").Value = "2.623.51"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("D19").Value = "'4.62"
$ws.Range("E19").Value = "  -2.60%  "
$ws.Range("D20").Value = "'340.96"
$ws.Range("E20").Value = "  -2.95%  "
$ws.Range("D21").Value = "'10.42"
$ws.Range("E21").Value = "  -1.53%  "
$ws.Range("D22").Value = "'6.11"
$ws.Range("E22").Value = "  -1.73%  "
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("D24").Value = "'60.65"
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("D25").Value = "'0.419"
$ws.Range("E25").Value = "  -2.06%  "
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").Value = "'0.161"
$ws.Range("E27").Value = "  -2.87%  "
$ws.Range("D28").Value = "0.0₃0806"
$ws.Range("E28").Value = "  -4.40%  "
$ws.Range("D29").Value = "'7.03"
$ws.Range("E29").Value = "  -4.18%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "'5.97"
$ws.Range("E32").Value = "  -4.97%  "
$ws.Range("D33").Value = "'18.95"
$ws.Range("E33").Value = "  -2.24%  "
$ws.Range("D34").Value = "'149.85"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "'3.94"
$ws.Range("E35").Value = "  -5.13%  "
$ws.Range("D36").Value = "'0.915"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "'1.13"
$ws.Range("E37").Value = "  -5.22%  "
$ws.Range("D38").Value = "'0.861"
$ws.Range("E38").Value = "  +2.46%  "
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("D40").Value = "'1.42"
$ws.Range("E40").Value = "  -4.34%  "
$ws.Range("D41").Value = "'3.62"
$ws.Range("E41").Value = "  -4.39%  "
$ws.Range("D42").Value = "'288.16"
$ws.Range("E42").Value = "  +0.90%  "
$ws.Range("D43").Value = "'0.624"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("E44").Value = "  -1.20%  "
$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").Value = "'0.0546"
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("D47").Value = "'19.43"
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("D48").Value = "'10.40"
$ws.Range("E48").Value = "  +0.92%  "
$ws.Range("D49").Value = "'0.0231"
$ws.Range("E49").Value = "  -2.36%  "
$ws.Range("D50").Value = "'4.68"
$ws.Range("E50").Value = "  -5.21%  "
$ws.Range("D51").Value = "1.963.49"
$ws.Range("E51").Value = "  -0.16%  "
